$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.417.58"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.877.04"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7142"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "239.00"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.9999"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07867"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.97%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3078"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "25.42"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +9.49%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08198"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.880.10"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.272"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.7249"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.99%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "89.48"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "29.478.98"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.835"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "242.41"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.58%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007847"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.35"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.132.16"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.9995"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.802"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.62%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1472"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +2.32%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "162.41"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.988"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.20"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.942"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.361"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -4.91%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.478"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.334"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.10%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.090"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05242"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.87%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.194"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.24%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7216"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.96%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +0.02%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01858"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.04%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.703"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.174.01"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9123"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.18%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "6.007"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.67%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4315"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "71.84"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.92%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "102.57"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5338"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.49%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.774"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.942"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +7.01%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "9.236"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
